# Update validation results for Wilke models.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the measured "Wilke [Mpa]" source values in column C.
# The dependent "Wilke %" formulas in column D (=Cn/C5) recalc automatically.
$ws.Range("C2").Value = 109.1407
$ws.Range("C3").Value = 266.221
$ws.Range("C4").Value = 482.8453
$ws.Range("C5").Value = 541.4263
$ws.Range("C6").Value = 1228.507
$ws.Range("C7").Value = 1126.4
$ws.Range("C8").Value = 2498.25
$ws.Range("C9").Value = 1961.558

# Match the saved cursor/selection position.
$ws.Range("H7").Select()
